# Applies the "Filled out Board and BoardCell classes..." edit to map.xlsx:
#   - Selects A1:Y24 (the populated board area) on the "map" sheet, with the
#     active cell anchored at the bottom-right corner of that range (Y24).
#   - Clears the helper "row index" formulas/values out of column Z (rows 1-24)
#     while leaving their cell formatting (style) untouched.
#   - Clears the helper "column index" formulas/values out of row 25
#     (columns A-Y) while leaving their cell formatting (style) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection -----------------------------------------------------------
# Select the used board range and anchor the active cell at its last cell
# (bottom-right corner), matching a drag/shift-click selection ending at Y24.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(24, 25))
$fullRange.Select()
$ws.Cells.Item(24, 25).Activate()
$fullRange.Select()

# --- Clear the old running-count helper formulas in column Z (rows 1-24) -
$ws.Range("Z1:Z24").ClearContents()

# --- Clear the old running-count helper formulas in row 25 (columns A-Y) -
$ws.Range("A25:Y25").ClearContents()
